$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name, Link) ---
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"

# --- Numeric-looking text columns (Price, Volume, Hora) ---
# Force text storage (matches the original inlineStr text cells): apply a
# Text number format to each contiguous block, write the literal values,
# then clear the formatting again so the cell style index reverts to the
# workbook default and no stray style is left behind. NOTE: this COM host
# only honours NumberFormat/ClearFormats on the FIRST area of a multi-area
# (comma-joined) Range, so every contiguous block gets its own statement
# instead of being combined with commas.
$ws.Range("D2:D15").NumberFormat = "@"
$ws.Range("D17:D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23:D25").NumberFormat = "@"
$ws.Range("D38:D45").NumberFormat = "@"
$ws.Range("D47:D50").NumberFormat = "@"
$ws.Range("E2:E26").NumberFormat = "@"
$ws.Range("E38:E50").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "329.32"
$ws.Range("E2").Value = "0.55%"
$ws.Range("G2").Value = "20"
$ws.Range("D3").Value = "44.32"
$ws.Range("E3").Value = "0.34%"
$ws.Range("G3").Value = "20"
$ws.Range("D4").Value = "5.570"
$ws.Range("E4").Value = "1.95%"
$ws.Range("G4").Value = "20"
$ws.Range("D5").Value = "0.08097"
$ws.Range("E5").Value = "0.22%"
$ws.Range("G5").Value = "20"
$ws.Range("D6").Value = "1.985"
$ws.Range("E6").Value = "5.14%"
$ws.Range("G6").Value = "20"
$ws.Range("D7").Value = "4.330"
$ws.Range("E7").Value = "0.78%"
$ws.Range("G7").Value = "20"
$ws.Range("D8").Value = "0.9538"
$ws.Range("E8").Value = "1.64%"
$ws.Range("G8").Value = "20"
$ws.Range("D9").Value = "2.569"
$ws.Range("E9").Value = "-5.04%"
$ws.Range("G9").Value = "20"
$ws.Range("D10").Value = "0.1174"
$ws.Range("E10").Value = "-1.91%"
$ws.Range("G10").Value = "20"
$ws.Range("D11").Value = "0.1854"
$ws.Range("E11").Value = "-2.00%"
$ws.Range("G11").Value = "20"
$ws.Range("D12").Value = "10.23"
$ws.Range("E12").Value = "19.31%"
$ws.Range("G12").Value = "20"
$ws.Range("D13").Value = "0.09936"
$ws.Range("E13").Value = "2.06%"
$ws.Range("G13").Value = "20"
$ws.Range("D14").Value = "0.04744"
$ws.Range("E14").Value = "16.33%"
$ws.Range("G14").Value = "20"
$ws.Range("D15").Value = "0.1069"
$ws.Range("E15").Value = "0.05%"
$ws.Range("G15").Value = "20"
$ws.Range("E16").Value = "0.94%"
$ws.Range("G16").Value = "20"
$ws.Range("D17").Value = "0.04226"
$ws.Range("E17").Value = "-2.93%"
$ws.Range("G17").Value = "20"
$ws.Range("D18").Value = "0.005903"
$ws.Range("E18").Value = "-1.27%"
$ws.Range("G18").Value = "20"
$ws.Range("E19").Value = "-5.76%"
$ws.Range("G19").Value = "20"
$ws.Range("E20").Value = "-0.76%"
$ws.Range("G20").Value = "20"
$ws.Range("D21").Value = "0.1410"
$ws.Range("E21").Value = "6.03%"
$ws.Range("G21").Value = "20"
$ws.Range("E22").Value = "0.41%"
$ws.Range("G22").Value = "20"
$ws.Range("D23").Value = "0.001250"
$ws.Range("E23").Value = "1.19%"
$ws.Range("G23").Value = "20"
$ws.Range("D24").Value = "0.004364"
$ws.Range("E24").Value = "2.06%"
$ws.Range("G24").Value = "20"
$ws.Range("D25").Value = "0.0001192"
$ws.Range("E25").Value = "-3.53%"
$ws.Range("G25").Value = "20"
$ws.Range("E26").Value = "-0.71%"
$ws.Range("G26").Value = "20"
$ws.Range("G27").Value = "20"
$ws.Range("G28").Value = "20"
$ws.Range("G29").Value = "20"
$ws.Range("G30").Value = "20"
$ws.Range("G31").Value = "20"
$ws.Range("G32").Value = "20"
$ws.Range("G33").Value = "20"
$ws.Range("G34").Value = "20"
$ws.Range("G35").Value = "20"
$ws.Range("G36").Value = "20"
$ws.Range("G37").Value = "20"
$ws.Range("D38").Value = "0.02651"
$ws.Range("E38").Value = "-0.15%"
$ws.Range("G38").Value = "20"
$ws.Range("D39").Value = "0.05554"
$ws.Range("E39").Value = "1.96%"
$ws.Range("G39").Value = "20"
$ws.Range("D40").Value = "0.007587"
$ws.Range("E40").Value = "-1.45%"
$ws.Range("G40").Value = "20"
$ws.Range("D41").Value = "0.1408"
$ws.Range("E41").Value = "1.49%"
$ws.Range("G41").Value = "20"
$ws.Range("D42").Value = "0.008078"
$ws.Range("E42").Value = "-16.92%"
$ws.Range("G42").Value = "20"
$ws.Range("D43").Value = "0.002019"
$ws.Range("E43").Value = "-4.85%"
$ws.Range("G43").Value = "20"
$ws.Range("D44").Value = "0.008907"
$ws.Range("E44").Value = "-10.03%"
$ws.Range("G44").Value = "20"
$ws.Range("D45").Value = "0.00007215"
$ws.Range("E45").Value = "1.19%"
$ws.Range("G45").Value = "20"
$ws.Range("E46").Value = "-0.29%"
$ws.Range("G46").Value = "20"
$ws.Range("D47").Value = "0.003520"
$ws.Range("E47").Value = "-1.66%"
$ws.Range("G47").Value = "20"
$ws.Range("D48").Value = "0.002272"
$ws.Range("E48").Value = "-0.36%"
$ws.Range("G48").Value = "20"
$ws.Range("D49").Value = "0.00002104"
$ws.Range("E49").Value = "-0.29%"
$ws.Range("G49").Value = "20"
$ws.Range("D50").Value = "0.0002003"
$ws.Range("E50").Value = "-0.29%"
$ws.Range("G50").Value = "20"
$ws.Range("G51").Value = "20"

$ws.Range("D2:D15").ClearFormats()
$ws.Range("D17:D18").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23:D25").ClearFormats()
$ws.Range("D38:D45").ClearFormats()
$ws.Range("D47:D50").ClearFormats()
$ws.Range("E2:E26").ClearFormats()
$ws.Range("E38:E50").ClearFormats()
$ws.Range("G2:G51").ClearFormats()
